$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.002, -0.451, 0.239),
    @(-0.75, -0.023, -0.046),
    @(-0.198, 0.509, -0.394),
    @(0.25, 0.694, 0.396),
    @(0.011, -0.003, 0.531),
    @(-0.253, -0.671, 0.035),
    @(-0.135, 0.363, -0.324),
    @(0.294, -0.301, 0.795),
    @(0.362, -0.036, 0.196),
    @(-0.266, -0.226, -0.196)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
